$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set values in the same order the original author entered them, so the
# generated shared-strings table matches the canonical ordering.
$ws.Range("A1").Value = "Nén mp3"
$ws.Range("B2").Value = "http://www.yoyogames.com/resources?cat_id=4"
$ws.Range("B3").Value = "The Witcher 2 Bonus Disc"
$ws.Range("B5").Value = "Bejeweled 2"
$ws.Range("B4").Value = "http://www.gamedev.net/topic/272386-sprites-sprites-and-more-sprites/"
$ws.Range("A2").Value = "Các tập tin âm thanh (mp3, wav)"
$ws.Range("A4").Value = "Các tập tin hình ảnh (jpg)"
$ws.Range("B1").Value = "EKOS MP3Minimizer"
$ws.Range("B6").Value = "Resource cung cấp trong môn học C4W"

$ws.Hyperlinks.Add($ws.Range("B2"), "http://www.yoyogames.com/resources?cat_id=4")
$ws.Hyperlinks.Add($ws.Range("B4"), "http://www.gamedev.net/topic/272386-sprites-sprites-and-more-sprites/")

# The host engine's ColumnWidth setter snaps to 1/6-character increments
# (raw = round(input*6)/6 + 5/6), so these inputs are chosen to land the
# stored <col> width attribute as close as possible to the canonical
# 29.85546875 / 68.5703125 (Excel's 1/256-character units) targets.
$ws.Columns.Item(1).ColumnWidth = 29.0
$ws.Columns.Item(2).ColumnWidth = 67.66666666666667
